$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new "Description" column is being inserted as column C. Columns C..G
# (Valor Accion, Locator, Valor Locator, Screenshot, Wait Time) all shift one
# place to the right, becoming D..H. The column-width metadata must stay put
# (it still targets the same physical columns as before), so we move the
# cell contents/formats by hand instead of doing a true Insert.
# ---------------------------------------------------------------------------

# Shift the data in C:G one column to the right (into D:H), working from the
# rightmost column first so nothing gets clobbered before it is read. Each
# destination is cleared before the copy so that blank source cells really
# blank out the destination (PasteSpecial of an empty cell is a no-op here).
# The shift is done with two PasteSpecial passes: xlPasteAll (-4104) for the
# values, then xlPasteFormats (-4122) to make sure number-format-driven
# style differences (e.g. the "@" text format) really land on the target,
# since -4104 alone doesn't carry NumberFormat in this host.
$ws.Range("H1:H7").ClearContents()
$ws.Range("G1:G7").Copy()
$ws.Range("H1").PasteSpecial(-4104)
$ws.Range("G1:G7").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("G1:G7").ClearContents()
$ws.Range("F1:F7").Copy()
$ws.Range("G1").PasteSpecial(-4104)
$ws.Range("F1:F7").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("F1:F7").ClearContents()
$ws.Range("E1:E7").Copy()
$ws.Range("F1").PasteSpecial(-4104)
$ws.Range("E1:E7").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("E1:E7").ClearContents()
$ws.Range("D1:D7").Copy()
$ws.Range("E1").PasteSpecial(-4104)
$ws.Range("D1:D7").Copy()
$ws.Range("E1").PasteSpecial(-4122)

$ws.Range("D1:D7").ClearContents()
$ws.Range("C1:C7").Copy()
$ws.Range("D1").PasteSpecial(-4104)
$ws.Range("C1:C7").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Move the hyperlink that used to live on C2 (Valor Accion) over to its new
# home on D2.
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.google.com/")
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# The old column C (now duplicated into D) still holds its original values;
# clear it out so the new Description column starts blank.
$ws.Range("C1:C7").ClearContents()

# Column H was freshly extended past the old used range, so borrow the
# neighbouring row styling for it (and for the one stray un-styled cell,
# E6, that the shift exposed).
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("H2:H7").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# Give the new column C the same banded look as the rest of the table.
$ws.Range("A2").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("C2:C7").PasteSpecial(-4122)

# Header renames + the new column's content.
$ws.Range("B1").Value = "Action"
$ws.Range("C1").Value = "Description"
$ws.Range("C2").Value = "Navegar a google"

$ws.Application.CutCopyMode = $false

# Match the recorded selection state.
$ws.Range("C9").Select()
